# Add a new service-history event ("حدث جديد") to Card14:
#  - fill the previously-blank row 16 cells (B:K, M) with the literal text "nan"
#  - append a new row 17 with the new event data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card14")

# --- Row 16: backfill empty cells with literal "nan" text (B:K, M) ---
$row16Cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "M")
foreach ($col in $row16Cols) {
    $ws.Range($col + "16").Value = "nan"
}

# --- Row 17: new event row ---
# A17 looks numeric ("14") but must be stored as text like the rest of the
# sheet, so force a text format, assign, then restore the default style so
# no stray number-format style is left behind on the cell.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "14"
$ws.Range("A17").Style = "Normal"

$ws.Range("L17").Value = "22\4\2025"
$ws.Range("M17").Value = "578.1 t"
$ws.Range("N17").Value = "تم تغيير الجرائد الاماميه (1_2_4_5_7_8) ومعيارته"
$ws.Range("O17").Value = "الخبير"
